$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = -4
    8  = -3
    10 = -6
    11 = -7
    12 = -5
    18 = 0
    21 = -1
    23 = -4
    24 = -3
    31 = -6
    32 = -2
    35 = -1
    38 = -6
    39 = -2
    42 = -4
    43 = -8
    48 = 14
    51 = 6
    54 = -9
    55 = 1
    58 = -5
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
